$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.197.45'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.904.27'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.80%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '306.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5254'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.96%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3778'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.53%  '
$ws.Range("E9").Value = '  +0.59%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.15'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8999'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.39%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08372'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +9.93%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.914.32'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.40%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.93'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.275'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.07%  '
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008607'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.55'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9999'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.236.35'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.41%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.067'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.22%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.141.96'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.28%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.61'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.438'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.20%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.47'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.23%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.285'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.83%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.752'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.30%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.17'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.75%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '114.74'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.926'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.30%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.817'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09289'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8090'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.67%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05064'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("E35").Value = '  +3.52%  '
$ws.Range("E36").Value = '  -2.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.371'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.84%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.619'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.17%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5718'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01989'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.14%  '
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.645'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.984'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.07%  '
$ws.Range("E44").Value = '  -0.31%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1512'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.40%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4844'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.20'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.54%  '
$ws.Range("E48").Value = '  -0.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.614'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '37.46'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.72%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.75'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.37%  '
